# Contactor Driver BOM update:
# - Change quantity of Molex Microfit 1x04 Vertical connector (row 6) from 4 to 1
# - Add a new BOM line (row 13) for "Molex Microfit Connector Horizontal"
# - Add a "Total" row (row 14) summing the cost column
# - Clear out the previously-blank placeholder rows (15-35) that no longer carry
#   stray C-column zero values / shared formulas
# - Update the active selection to D15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6: quantity update (Molex Microfit 1x04 Vertical) ---
$ws.Range("B6").Value2 = 1

# --- Row 13: new BOM line for the horizontal connector ---
# (Leave D13's inherited shared formula untouched so the shared formula group
#  for column D is preserved as long as possible.)
$ws.Range("A13").Value2 = "Molex Microfit Connector Horizontal"
$ws.Range("B13").Value2 = 3
$ws.Range("C13").Value2 = 0.39
$ws.Range("E13").Value2 = "https://www.mouser.com/ProductDetail/Molex/43645-0200?qs=4XSMV6Twtb2TZ7elJDViLw%3D%3D&gclid=CjwKCAjwpayjBhAnEiwA-7enaz5rrXZ-_uL_IGfOE-nYsLismrqezJIPpDvBGMZQ0Ird-ldIHxlpyhoCJCwQAvD_BwE"
$ws.Range("F13").Value2 = "https://www.mouser.com/datasheet/2/276/3/0436450200_CRIMP_HOUSINGS-2866614.pdf"

# --- Row 14: Total row ---
$ws.Range("C1").Copy()
[void]$ws.Range("C14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C14").Value2 = "Total"
$ws.Range("D14").Formula = "=SUM(D2:D13)"
$ws.Rows(14).RowHeight = 15.75

# --- Rows 15-35: clear the placeholder zero values/formulas ---
# Column C: fully clear (contents + formatting) so the cells disappear entirely
$ws.Range("C15:C35").Clear()
# Column D: keep the existing currency style, only remove the formula/value
$ws.Range("D15:D35").ClearContents()

# --- Update selection to D15 to match the saved workbook state ---
[void]$ws.Range("D15").Select()

Write-Output "Edit complete"
